$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Stable "donor" cells elsewhere on the board that already carry the exact
# fill style we need, and are themselves untouched by this edit. Copying
# their format (instead of poking Interior.Color directly) makes the saved
# workbook reuse the existing theme-based fill/style entries rather than
# inventing new literal-RGB ones.
$plainDonor = $ws.Range("B2")   # style 1 - plain room-letter fill
$doorDonor  = $ws.Range("D5")   # style 6 - door-cell fill ("#" markers)
$starDonor  = $ws.Range("L5")   # style 8 - room-center fill used by "S*"/R10

function Copy-ClueStyle($donor, $addr) {
    $donor.Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}

# Room H (top-left) - move the door marker from D4 to C4
Copy-ClueStyle $doorDonor  "C4"
Copy-ClueStyle $plainDonor "D4"
$ws.Range("C4").Value2 = "H#"
$ws.Range("D4").Value2 = "H"

# Room C (top-middle) - move the door marker from L4 to K4
Copy-ClueStyle $doorDonor  "K4"
Copy-ClueStyle $plainDonor "L4"
$ws.Range("K4").Value2 = "C#"
$ws.Range("L4").Value2 = "C"

# Room P (top-right) - move the door marker from T4 to S4
Copy-ClueStyle $doorDonor  "S4"
Copy-ClueStyle $plainDonor "T4"
$ws.Range("S4").Value2 = "P#"
$ws.Range("T4").Value2 = "P"

# Room S (left) - rotate door/plain/center markers between E10, F10, R10
Copy-ClueStyle $doorDonor  "E10"
Copy-ClueStyle $plainDonor "F10"
Copy-ClueStyle $starDonor  "R10"
$ws.Range("E10").Value2 = "S#"
$ws.Range("F10").Value2 = "S"
$ws.Range("R10").Value2 = "S*"

# Room D (bottom-left) - move the door marker from D15 to B15
Copy-ClueStyle $doorDonor  "B15"
Copy-ClueStyle $plainDonor "D15"
$ws.Range("B15").Value2 = "D#"
$ws.Range("D15").Value2 = "D"

# Room L (bottom-middle) - move the door marker from L15 to J15
Copy-ClueStyle $doorDonor  "J15"
Copy-ClueStyle $plainDonor "L15"
$ws.Range("J15").Value2 = "L#"
$ws.Range("L15").Value2 = "L"

# Room G (bottom-right) - move the door marker from T15 to R15
Copy-ClueStyle $doorDonor  "R15"
Copy-ClueStyle $plainDonor "T15"
$ws.Range("R15").Value2 = "G#"
$ws.Range("T15").Value2 = "G"

# Room M (bottom-left, row 21) - move the door marker C21<-, center D21<-,
# and move the old door marker (J21) back to a plain cell while K21 becomes
# the new room-center cell
Copy-ClueStyle $doorDonor  "C21"
Copy-ClueStyle $plainDonor "D21"
Copy-ClueStyle $plainDonor "J21"
Copy-ClueStyle $doorDonor  "K21"
$ws.Range("C21").Value2 = "M#"
$ws.Range("D21").Value2 = "M"
$ws.Range("J21").Value2 = "M"
$ws.Range("K21").Value2 = "M*"

# Room E (bottom-right, row 21) - move the door marker from U21 to O21
Copy-ClueStyle $doorDonor  "O21"
Copy-ClueStyle $plainDonor "U21"
$ws.Range("O21").Value2 = "E#"
$ws.Range("U21").Value2 = "E"

# Update the active selection to match the saved workbook state
$ws.Range("X21").Select()
